$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (raw OOXML width = ColumnWidth + 0.8333333333333333) ---
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667  # column B -> raw width 8
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667  # column G -> raw width 8
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667  # column K -> raw width 8
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667  # column L -> raw width 8
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667  # column M -> raw width 8
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667  # column O -> raw width 8
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667  # column P -> raw width 8
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666  # column T -> raw width 9
$ws.Columns.Item(33).ColumnWidth = 6.166666666666667  # column AG -> raw width 7
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667  # column AH -> raw width 8

# --- Update data rows 2-5 with new values ---
$rowVals = @(45124.50694444445, 4.96, 4.926, 1.074, 9.69, 9.442, 3.413, 8.191000000000001, 4.526, 2.102, 4.543, 5.171, 3.55, 0.877, 3.454, 4.766, 2.485, 1.205, 0, 48.94, 9.632, 2.995, 5.545, 4.641, 0.661, 7.552, 1.905, 4.889, 2.236, 4.485, 0.29, 5.812, 2.376, 3.752)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $rowVals[$c - 1]
}

$rowVals = @(45124.51388888889, 0.377, 0.784, 0.336, 0.6879999999999999, 0.765, 0, 3.092, 0.238, 0.173, 0.731, 0.447, 0, 0, 0.11, 0.544, 0.173, 0.432, 0, 0, 1.052, 0.07199999999999999, 0.7, 0.929, 0.073, 3.018, 0.111, 1.327, 0.107, 0.575, 0.034, 2.898, 0.273, 0.143)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(3, $c).Value = $rowVals[$c - 1]
}

$rowVals = @(45124.52083333334, 15.356, 11.81, 0.73, 33.212, 27.579, 12.273, 39.817, 18.383, 8.241, 12.717, 13.429, 13.934, 3.76, 11.838, 17.081, 9.906000000000001, 0.419, 0.421, 176.219, 33.156, 10.91, 22.374, 12.19, 1.586, 21.212, 9.640000000000001, 9.282, 9.957000000000001, 14.22, 0, 35.664, 6.416, 13.671)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(4, $c).Value = $rowVals[$c - 1]
}

$rowVals = @(45124.52777777778, 5.78, 4.52, 0.36, 12.37, 10.31, 4.74, 20.61, 6.81, 3.12, 4.85, 5.04, 5.14, 1.38, 4.36, 6.46, 3.7, 0.27, 0.13, 61.65, 12.59, 4.01, 8.52, 4.7, 0.6, 10.53, 3.53, 3.73, 3.66, 5.38, 0, 19.01, 2.45, 5.05)
for ($c = 1; $c -le $rowVals.Length; $c++) {
    $ws.Cells.Item(5, $c).Value = $rowVals[$c - 1]
}

# --- Remove the now-obsolete last data row (row 6) ---
$ws.Rows.Item(6).Delete()
